$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A to hold the "Match ID" field.
$ws.Columns("A").Insert()

# Header row (row 1): new "Match ID" header, bold like the other headers
# but without the border/center-alignment used on the rest of row 1.
$ws.Range("A1").Value = "Match ID"
$ws.Range("A1").Font.Bold = $true

# Row 2 is a hidden duplicate header row; its new A2 cell stays blank but
# picks up the same bold style, and its old "Player" label (now in B2)
# is cleared out.
$ws.Range("A2").Font.Bold = $true
$ws.Range("B2").ClearContents()

# Row 3 is a hidden spacer row; give its new A3 cell the same style too.
$ws.Range("A3").Font.Bold = $true

# Data + totals rows (4-20): fill in the Match ID value.
$ws.Range("A4:A20").Value = 31
$ws.Range("A4:A20").Font.Bold = $true

# Row 20 (the hidden totals row) picks up a stray custom row-height when its
# cell value is written; auto-fit it back to the sheet default so the row
# stays attribute-for-attribute the same as before the edit.
$ws.Rows(20).AutoFit()

# Reflect the new selection left in the sheet after the edit.
$ws.Range("A4:A20").Select() | Out-Null
